$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5 (pushes existing rows 5..41 down to 6..42),
# copying formatting from the row above (matches Excel's default Insert
# behaviour, which also keeps the date number format on column D).
$ws.Range("A5").EntireRow.Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44819
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100114007
$ws.Range("G5").Value = "Jengibre"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 13400
$ws.Range("N5").Value = '$/caja 13 kilos'
$ws.Range("O5").Value = "Perú"
$ws.Range("P5").Value = 1031
$ws.Range("Q5").Value = 13
$ws.Range("R5").Value = "Hortaliza"
